# Add a new "sub_module" row to the "summary" sheet, describing the new
# sub_module object added to verigen, and make "summary" the active sheet.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("summary")

# New row 5: A5="sub_module", B5="data" (centered style, like other B cells),
# C5="Sub-module instance list"
$summary.Cells.Item(5, 1).Value = "sub_module"
$summary.Cells.Item(5, 2).Value = "data"
$summary.Cells.Item(5, 2).HorizontalAlignment = -4108  # xlCenter
$summary.Cells.Item(5, 3).Value = "Sub-module instance list"

# Select the newly added cell and make "summary" the active sheet/tab.
$summary.Activate()
$summary.Range("C5").Select()
